$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 6299.75  # H74: 5812.125 -> 6299.75
$ws.Cells.Item(74, 9).Value = 6299.75  # I74: 5812.125 -> 6299.75
$ws.Cells.Item(74, 11).Value = 6299.75  # K74: 5812.125 -> 6299.75
$ws.Cells.Item(74, 13).Value = -5363.75  # M74: -4876.125 -> -5363.75
$ws.Cells.Item(76, 8).Value = 6991.75  # H76: 6883 -> 6991.75
$ws.Cells.Item(76, 9).Value = 6822.3335  # I76: 6759.6 -> 6822.3335
$ws.Cells.Item(76, 11).Value = 6822.3335  # K76: 6759.6 -> 6822.3335
$ws.Cells.Item(76, 13).Value = -6507.3335  # M76: -6444.6 -> -6507.3335
$ws.Cells.Item(77, 8).Value = 6299.75  # H77: 5812.125 -> 6299.75
$ws.Cells.Item(77, 9).Value = 6299.75  # I77: 5812.125 -> 6299.75
$ws.Cells.Item(77, 11).Value = 31498.75  # K77: 29060.625 -> 31498.75
$ws.Cells.Item(77, 13).Value = -26818.75  # M77: -24380.625 -> -26818.75
$ws.Cells.Item(79, 8).Value = 6991.75  # H79: 6883 -> 6991.75
$ws.Cells.Item(79, 9).Value = 6822.3335  # I79: 6759.6 -> 6822.3335
$ws.Cells.Item(79, 11).Value = 6822.3335  # K79: 6759.6 -> 6822.3335
$ws.Cells.Item(79, 13).Value = -5730.3335  # M79: -5667.6 -> -5730.3335
$ws.Cells.Item(92, 8).Value = 509.2857  # H92: 462.5 -> 509.2857
$ws.Cells.Item(92, 9).Value = 526.8333  # I92: 496.69232 -> 526.8333
$ws.Cells.Item(92, 10).Value = 404  # J92: 314.33334 -> 404
$ws.Cells.Item(92, 11).Value = 526.8333  # K92: 496.69232 -> 526.8333
$ws.Cells.Item(92, 12).Value = 404  # L92: 314.33334 -> 404
$ws.Cells.Item(92, 13).Value = 721.1667  # M92: 751.30768 -> 721.1667
$ws.Cells.Item(92, 14).Value = -2900  # N92: -2810.33334 -> -2900
$ws.Cells.Item(101, 8).Value = 1121.1428  # H101: 1031.125 -> 1121.1428
$ws.Cells.Item(101, 9).Value = 983.6667  # I101: 837.75 -> 983.6667
$ws.Cells.Item(101, 10).Value = 1224.25  # J101: 1224.5 -> 1224.25
$ws.Cells.Item(101, 11).Value = 2951.0001  # K101: 2513.25 -> 2951.0001
$ws.Cells.Item(101, 12).Value = 3672.75  # L101: 3673.5 -> 3672.75
$ws.Cells.Item(101, 13).Value = -1329.0001  # M101: -891.25 -> -1329.0001
$ws.Cells.Item(101, 14).Value = -6916.75  # N101: -6917.5 -> -6916.75
$ws.Cells.Item(125, 8).Value = 38463828  # H125: 38463550 -> 38463828
$ws.Cells.Item(125, 10).Value = 1742.4  # J125: 1032.4 -> 1742.4
$ws.Cells.Item(125, 12).Value = 15681.6  # L125: 9291.6 -> 15681.6
$ws.Cells.Item(125, 14).Value = -20601.6  # N125: -14211.6 -> -20601.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8200.588  # H32: 9455.357 -> 8200.588
$ws.Cells.Item(32, 9).Value = 4760.2334  # I32: 5731.2915 -> 4760.2334
$ws.Cells.Item(32, 10).Value = 34003.25  # J32: 31799.75 -> 34003.25
$ws.Cells.Item(32, 11).Value = 4760.2334  # K32: 5731.2915 -> 4760.2334
$ws.Cells.Item(32, 12).Value = 34003.25  # L32: 31799.75 -> 34003.25
$ws.Cells.Item(32, 13).Value = -4473.2334  # M32: -5444.2915 -> -4473.2334
$ws.Cells.Item(32, 14).Value = -34577.25  # N32: -32373.75 -> -34577.25
$ws.Cells.Item(34, 8).Value = 50555.332  # H34: 51249.75 -> 50555.332
$ws.Cells.Item(34, 10).Value = 56428.43  # J34: 58333.168 -> 56428.43
$ws.Cells.Item(34, 12).Value = 56428.43  # L34: 58333.168 -> 56428.43
$ws.Cells.Item(34, 14).Value = -56970.43  # N34: -58875.168 -> -56970.43
$ws.Cells.Item(63, 8).Value = 5302.7715  # H63: 5070.8125 -> 5302.7715
$ws.Cells.Item(63, 10).Value = 8166.706  # J63: 8250.214 -> 8166.706
$ws.Cells.Item(63, 12).Value = 8166.706  # L63: 8250.214 -> 8166.706
$ws.Cells.Item(63, 14).Value = -9538.706  # N63: -9622.214 -> -9538.706
$ws.Cells.Item(66, 8).Value = 5302.7715  # H66: 5070.8125 -> 5302.7715
$ws.Cells.Item(66, 10).Value = 8166.706  # J66: 8250.214 -> 8166.706
$ws.Cells.Item(66, 12).Value = 40833.53  # L66: 41251.07 -> 40833.53
$ws.Cells.Item(66, 14).Value = -47697.53  # N66: -48115.07 -> -47697.53
$ws.Cells.Item(88, 8).Value = 3279.8  # H88: 3499.6667 -> 3279.8
$ws.Cells.Item(88, 9).Value = 2900  # I88: 0 -> 2900
$ws.Cells.Item(88, 10).Value = 3374.75  # J88: 3499.6667 -> 3374.75
$ws.Cells.Item(88, 11).Value = 2900  # K88: 0 -> 2900
$ws.Cells.Item(88, 12).Value = 3374.75  # L88: 3499.6667 -> 3374.75
$ws.Cells.Item(88, 13).Value = -2494  # M88: None -> -2494
$ws.Cells.Item(88, 14).Value = -4186.75  # N88: -4311.6667 -> -4186.75
$ws.Cells.Item(91, 8).Value = 3279.8  # H91: 3499.6667 -> 3279.8
$ws.Cells.Item(91, 9).Value = 2900  # I91: 0 -> 2900
$ws.Cells.Item(91, 10).Value = 3374.75  # J91: 3499.6667 -> 3374.75
$ws.Cells.Item(91, 11).Value = 2900  # K91: 0 -> 2900
$ws.Cells.Item(91, 12).Value = 3374.75  # L91: 3499.6667 -> 3374.75
$ws.Cells.Item(91, 13).Value = -1496  # M91: None -> -1496
$ws.Cells.Item(91, 14).Value = -6182.75  # N91: -6307.6667 -> -6182.75
$ws.Cells.Item(132, 8).Value = 5442.149  # H132: 4747.14 -> 5442.149
$ws.Cells.Item(132, 9).Value = 3784.1052  # I132: 3304.2083 -> 3784.1052
$ws.Cells.Item(132, 11).Value = 11352.3156  # K132: 9912.624899999999 -> 11352.3156
$ws.Cells.Item(132, 13).Value = -8822.3156  # M132: -7382.624899999999 -> -8822.3156

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 4118.2  # H12: 2648 -> 4118.2
$ws.Cells.Item(12, 10).Value = 9999  # J12: 0 -> 9999
$ws.Cells.Item(12, 12).Value = 9999  # L12: 0 -> 9999
$ws.Cells.Item(12, 14).Value = -10335  # N12: None -> -10335
$ws.Cells.Item(20, 8).Value = 3767.7144  # H20: 2612.5386 -> 3767.7144
$ws.Cells.Item(20, 9).Value = 4559.6665  # I20: 2283.75 -> 4559.6665
$ws.Cells.Item(20, 10).Value = 3173.75  # J20: 3138.6 -> 3173.75
$ws.Cells.Item(20, 11).Value = 4559.6665  # K20: 2283.75 -> 4559.6665
$ws.Cells.Item(20, 12).Value = 3173.75  # L20: 3138.6 -> 3173.75
$ws.Cells.Item(20, 13).Value = -4312.6665  # M20: -2036.75 -> -4312.6665
$ws.Cells.Item(20, 14).Value = -3667.75  # N20: -3632.6 -> -3667.75
$ws.Cells.Item(86, 8).Value = 670298.3  # H86: 457681.6 -> 670298.3
$ws.Cells.Item(86, 9).Value = 1431756.8  # I86: 771870.6 -> 1431756.8
$ws.Cells.Item(86, 10).Value = 4022.25  # J86: 3853 -> 4022.25
$ws.Cells.Item(86, 11).Value = 1431756.8  # K86: 771870.6 -> 1431756.8
$ws.Cells.Item(86, 12).Value = 4022.25  # L86: 3853 -> 4022.25
$ws.Cells.Item(86, 13).Value = -1430633.8  # M86: -770747.6 -> -1430633.8
$ws.Cells.Item(86, 14).Value = -6268.25  # N86: -6099 -> -6268.25
$ws.Cells.Item(89, 8).Value = 670298.3  # H89: 457681.6 -> 670298.3
$ws.Cells.Item(89, 9).Value = 1431756.8  # I89: 771870.6 -> 1431756.8
$ws.Cells.Item(89, 10).Value = 4022.25  # J89: 3853 -> 4022.25
$ws.Cells.Item(89, 11).Value = 7158784  # K89: 3859353 -> 7158784
$ws.Cells.Item(89, 12).Value = 20111.25  # L89: 19265 -> 20111.25
$ws.Cells.Item(89, 13).Value = -7153168  # M89: -3853737 -> -7153168
$ws.Cells.Item(89, 14).Value = -31343.25  # N89: -30497 -> -31343.25
$ws.Cells.Item(105, 8).Value = 9465.5  # H105: 4289.65 -> 9465.5
$ws.Cells.Item(105, 9).Value = 9249  # I105: 3812.25 -> 9249
$ws.Cells.Item(105, 10).Value = 9898.5  # J105: 6199.25 -> 9898.5
$ws.Cells.Item(105, 11).Value = 9249  # K105: 3812.25 -> 9249
$ws.Cells.Item(105, 12).Value = 9898.5  # L105: 6199.25 -> 9898.5
$ws.Cells.Item(105, 13).Value = -7502  # M105: -2065.25 -> -7502
$ws.Cells.Item(105, 14).Value = -13392.5  # N105: -9693.25 -> -13392.5
$ws.Cells.Item(107, 8).Value = 2937.5  # H107: 1987.4546 -> 2937.5
$ws.Cells.Item(107, 9).Value = 2125  # I107: 1750 -> 2125
$ws.Cells.Item(107, 10).Value = 3750  # J107: 2076.5 -> 3750
$ws.Cells.Item(107, 11).Value = 2125  # K107: 1750 -> 2125
$ws.Cells.Item(107, 12).Value = 3750  # L107: 2076.5 -> 3750
$ws.Cells.Item(107, 13).Value = -205  # M107: 170 -> -205
$ws.Cells.Item(107, 14).Value = -7590  # N107: -5916.5 -> -7590

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 3257  # H8: 3324.5 -> 3257
$ws.Cells.Item(8, 10).Value = 3535.2856  # J8: 3612.4285 -> 3535.2856
$ws.Cells.Item(8, 12).Value = 3535.2856  # L8: 3612.4285 -> 3535.2856
$ws.Cells.Item(8, 14).Value = -3815.2856  # N8: -3892.4285 -> -3815.2856
$ws.Cells.Item(31, 8).Value = 58892.688  # H31: 97579.91 -> 58892.688
$ws.Cells.Item(31, 9).Value = 53592.156  # I31: 126224.375 -> 53592.156
$ws.Cells.Item(31, 10).Value = 65187.062  # J31: 79952.53999999999 -> 65187.062
$ws.Cells.Item(31, 11).Value = 53592.156  # K31: 126224.375 -> 53592.156
$ws.Cells.Item(31, 12).Value = 65187.062  # L31: 79952.53999999999 -> 65187.062
$ws.Cells.Item(31, 13).Value = -53297.156  # M31: -125929.375 -> -53297.156
$ws.Cells.Item(31, 14).Value = -65777.06200000001  # N31: -80542.53999999999 -> -65777.06200000001
$ws.Cells.Item(33, 8).Value = 549  # H33: 10292.5 -> 549
$ws.Cells.Item(33, 9).Value = 549  # I33: 550 -> 549
$ws.Cells.Item(33, 10).Value = 0  # J33: 20035 -> 0
$ws.Cells.Item(33, 11).Value = 549  # K33: 550 -> 549
$ws.Cells.Item(33, 12).Value = 0  # L33: 20035 -> 0
$ws.Cells.Item(33, 13).Value = -170  # M33: -171 -> -170
$ws.Cells.Item(33, 14).ClearContents()  # N33: was -20793
$ws.Cells.Item(34, 8).Value = 58892.688  # H34: 97579.91 -> 58892.688
$ws.Cells.Item(34, 9).Value = 53592.156  # I34: 126224.375 -> 53592.156
$ws.Cells.Item(34, 10).Value = 65187.062  # J34: 79952.53999999999 -> 65187.062
$ws.Cells.Item(34, 11).Value = 53592.156  # K34: 126224.375 -> 53592.156
$ws.Cells.Item(34, 12).Value = 65187.062  # L34: 79952.53999999999 -> 65187.062
$ws.Cells.Item(34, 13).Value = -53390.156  # M34: -126022.375 -> -53390.156
$ws.Cells.Item(34, 14).Value = -65591.06200000001  # N34: -80356.53999999999 -> -65591.06200000001
$ws.Cells.Item(35, 8).Value = 2206.5454  # H35: 2388.182 -> 2206.5454
$ws.Cells.Item(35, 10).Value = 3994.5  # J35: 4993.5 -> 3994.5
$ws.Cells.Item(35, 12).Value = 3994.5  # L35: 4993.5 -> 3994.5
$ws.Cells.Item(35, 14).Value = -4582.5  # N35: -5581.5 -> -4582.5
$ws.Cells.Item(62, 8).Value = 503584.16  # H62: 603100.6 -> 503584.16
$ws.Cells.Item(62, 10).Value = 503501.5  # J62: 669334.7 -> 503501.5
$ws.Cells.Item(62, 12).Value = 503501.5  # L62: 669334.7 -> 503501.5
$ws.Cells.Item(62, 14).Value = -504749.5  # N62: -670582.7 -> -504749.5
$ws.Cells.Item(65, 8).Value = 503584.16  # H65: 603100.6 -> 503584.16
$ws.Cells.Item(65, 10).Value = 503501.5  # J65: 669334.7 -> 503501.5
$ws.Cells.Item(65, 12).Value = 2517507.5  # L65: 3346673.5 -> 2517507.5
$ws.Cells.Item(65, 14).Value = -2523747.5  # N65: -3352913.5 -> -2523747.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 3201.8  # H131: 3035.6956 -> 3201.8
$ws.Cells.Item(131, 9).Value = 605.2857  # I131: 623.1667 -> 605.2857
$ws.Cells.Item(131, 10).Value = 4599.923  # J131: 3887.1765 -> 4599.923
$ws.Cells.Item(131, 11).Value = 1815.8571  # K131: 1869.5001 -> 1815.8571
$ws.Cells.Item(131, 12).Value = 13799.769  # L131: 11661.5295 -> 13799.769
$ws.Cells.Item(131, 13).Value = 3224.1429  # M131: 3170.4999 -> 3224.1429
$ws.Cells.Item(131, 14).Value = -23879.769  # N131: -21741.5295 -> -23879.769
$ws.Cells.Item(132, 8).Value = 8335014.5  # H132: 7694011.5 -> 8335014.5
$ws.Cells.Item(132, 9).Value = 1767.7142  # I132: 1742.5 -> 1767.7142
$ws.Cells.Item(132, 10).Value = 20001560  # J132: 14287385 -> 20001560
$ws.Cells.Item(132, 11).Value = 15909.4278  # K132: 15682.5 -> 15909.4278
$ws.Cells.Item(132, 12).Value = 180014040  # L132: 128586465 -> 180014040
$ws.Cells.Item(132, 13).Value = -13379.4278  # M132: -13152.5 -> -13379.4278
$ws.Cells.Item(132, 14).Value = -180019100  # N132: -128591525 -> -180019100

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 2172.2856  # H9: 2315 -> 2172.2856
$ws.Cells.Item(9, 10).Value = 7000  # J9: 7499.5 -> 7000
$ws.Cells.Item(9, 12).Value = 7000  # L9: 7499.5 -> 7000
$ws.Cells.Item(9, 14).Value = -7340  # N9: -7839.5 -> -7340
$ws.Cells.Item(12, 8).Value = 15149.833  # H12: 15842.714 -> 15149.833
$ws.Cells.Item(12, 10).Value = 20499.5  # J12: 20333 -> 20499.5
$ws.Cells.Item(12, 12).Value = 20499.5  # L12: 20333 -> 20499.5
$ws.Cells.Item(12, 14).Value = -20779.5  # N12: -20613 -> -20779.5
$ws.Cells.Item(40, 8).Value = 38600  # H40: 40000 -> 38600
$ws.Cells.Item(40, 10).Value = 38600  # J40: 40000 -> 38600
$ws.Cells.Item(40, 12).Value = 38600  # L40: 40000 -> 38600
$ws.Cells.Item(40, 14).Value = -38902  # N40: -40302 -> -38902
$ws.Cells.Item(70, 8).Value = 11191.267  # H70: 11348.143 -> 11191.267
$ws.Cells.Item(70, 10).Value = 11485.625  # J70: 11841.429 -> 11485.625
$ws.Cells.Item(70, 12).Value = 11485.625  # L70: 11841.429 -> 11485.625
$ws.Cells.Item(70, 14).Value = -12025.625  # N70: -12381.429 -> -12025.625
$ws.Cells.Item(73, 8).Value = 11191.267  # H73: 11348.143 -> 11191.267
$ws.Cells.Item(73, 10).Value = 11485.625  # J73: 11841.429 -> 11485.625
$ws.Cells.Item(73, 12).Value = 11485.625  # L73: 11841.429 -> 11485.625
$ws.Cells.Item(73, 14).Value = -13357.625  # N73: -13713.429 -> -13357.625
$ws.Cells.Item(80, 8).Value = 3000  # H80: 0 -> 3000
$ws.Cells.Item(80, 10).Value = 3000  # J80: 0 -> 3000
$ws.Cells.Item(80, 12).Value = 3000  # L80: 0 -> 3000
$ws.Cells.Item(80, 14).Value = -4996  # N80: None -> -4996
$ws.Cells.Item(83, 8).Value = 3000  # H83: 0 -> 3000
$ws.Cells.Item(83, 10).Value = 3000  # J83: 0 -> 3000
$ws.Cells.Item(83, 12).Value = 15000  # L83: 0 -> 15000
$ws.Cells.Item(83, 14).Value = -24984  # N83: None -> -24984
$ws.Cells.Item(102, 8).Value = 2614  # H102: 2618.7368 -> 2614
$ws.Cells.Item(102, 9).Value = 2433.5715  # I102: 2481.5386 -> 2433.5715
$ws.Cells.Item(102, 10).Value = 3245.5  # J102: 2916 -> 3245.5
$ws.Cells.Item(102, 11).Value = 2433.5715  # K102: 2481.5386 -> 2433.5715
$ws.Cells.Item(102, 12).Value = 3245.5  # L102: 2916 -> 3245.5
$ws.Cells.Item(102, 13).Value = -811.5715  # M102: -859.5385999999999 -> -811.5715
$ws.Cells.Item(102, 14).Value = -6489.5  # N102: -6160 -> -6489.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4149  # H40: 4443.3335 -> 4149
$ws.Cells.Item(40, 9).Value = 4285.4287  # I40: 4749.6665 -> 4285.4287
$ws.Cells.Item(40, 11).Value = 4285.4287  # K40: 4749.6665 -> 4285.4287
$ws.Cells.Item(40, 13).Value = -4149.4287  # M40: -4613.6665 -> -4149.4287
$ws.Cells.Item(46, 8).Value = 1534.7273  # H46: 1564.6666 -> 1534.7273
$ws.Cells.Item(46, 9).Value = 1338.1111  # I46: 1528.6 -> 1338.1111
$ws.Cells.Item(46, 10).Value = 1608.4584  # J46: 1572.8636 -> 1608.4584
$ws.Cells.Item(46, 11).Value = 1338.1111  # K46: 1528.6 -> 1338.1111
$ws.Cells.Item(46, 12).Value = 1608.4584  # L46: 1572.8636 -> 1608.4584
$ws.Cells.Item(46, 13).Value = -1150.1111  # M46: -1340.6 -> -1150.1111
$ws.Cells.Item(46, 14).Value = -1984.4584  # N46: -1948.8636 -> -1984.4584
$ws.Cells.Item(68, 8).Value = 2681.125  # H68: 2603.0789 -> 2681.125
$ws.Cells.Item(68, 9).Value = 2414.0908  # I68: 2252.6453 -> 2414.0908
$ws.Cells.Item(68, 10).Value = 3940  # J68: 4155 -> 3940
$ws.Cells.Item(68, 11).Value = 2414.0908  # K68: 2252.6453 -> 2414.0908
$ws.Cells.Item(68, 12).Value = 3940  # L68: 4155 -> 3940
$ws.Cells.Item(68, 13).Value = -1665.0908  # M68: -1503.6453 -> -1665.0908
$ws.Cells.Item(68, 14).Value = -5438  # N68: -5653 -> -5438
$ws.Cells.Item(71, 8).Value = 2681.125  # H71: 2603.0789 -> 2681.125
$ws.Cells.Item(71, 9).Value = 2414.0908  # I71: 2252.6453 -> 2414.0908
$ws.Cells.Item(71, 10).Value = 3940  # J71: 4155 -> 3940
$ws.Cells.Item(71, 11).Value = 12070.454  # K71: 11263.2265 -> 12070.454
$ws.Cells.Item(71, 12).Value = 19700  # L71: 20775 -> 19700
$ws.Cells.Item(71, 13).Value = -8326.454  # M71: -7519.226500000001 -> -8326.454
$ws.Cells.Item(71, 14).Value = -27188  # N71: -28263 -> -27188

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 6000.143  # H17: 5211.778 -> 6000.143
$ws.Cells.Item(17, 9).Value = 5333.8335  # I17: 4613.5 -> 5333.8335
$ws.Cells.Item(17, 11).Value = 5333.8335  # K17: 4613.5 -> 5333.8335
$ws.Cells.Item(17, 13).Value = -5161.8335  # M17: -4441.5 -> -5161.8335
$ws.Cells.Item(32, 8).Value = 4242  # H32: 4408.6665 -> 4242
$ws.Cells.Item(32, 10).Value = 4350  # J32: 4600 -> 4350
$ws.Cells.Item(32, 12).Value = 4350  # L32: 4600 -> 4350
$ws.Cells.Item(32, 14).Value = -4984  # N32: -5234 -> -4984
$ws.Cells.Item(42, 8).Value = 66011  # H42: 66011.5 -> 66011
$ws.Cells.Item(42, 10).Value = 85023  # J42: 85024 -> 85023
$ws.Cells.Item(42, 12).Value = 85023  # L42: 85024 -> 85023
$ws.Cells.Item(42, 14).Value = -85779  # N42: -85780 -> -85779
$ws.Cells.Item(43, 8).Value = 43124.5  # H43: 43124.75 -> 43124.5
$ws.Cells.Item(43, 10).Value = 38999.5  # J43: 39000 -> 38999.5
$ws.Cells.Item(43, 12).Value = 38999.5  # L43: 39000 -> 38999.5
$ws.Cells.Item(43, 14).Value = -39297.5  # N43: -39298 -> -39297.5
$ws.Cells.Item(100, 8).Value = 222.46666  # H100: 202.61111 -> 222.46666
$ws.Cells.Item(100, 9).Value = 135.45454  # I100: 128.57143 -> 135.45454
$ws.Cells.Item(100, 11).Value = 270.90908  # K100: 257.14286 -> 270.90908
$ws.Cells.Item(100, 13).Value = 270.09092  # M100: 283.85714 -> 270.09092
